# Update "想去人数" (interest count) values in column F for both the
# "展览" sheet and the "全部类型" sheet, reflecting refreshed counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4392
$wsExhibit.Range("F3").Value = 2463
$wsExhibit.Range("F5").Value = 25
$wsExhibit.Range("F8").Value = 221
$wsExhibit.Range("F11").Value = 163
$wsExhibit.Range("F12").Value = 1648
$wsExhibit.Range("F14").Value = 3517
$wsExhibit.Range("F15").Value = 233

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4392
$wsAll.Range("F3").Value = 2463
$wsAll.Range("F5").Value = 25
$wsAll.Range("F10").Value = 221
$wsAll.Range("F13").Value = 163
$wsAll.Range("F16").Value = 1648
$wsAll.Range("F18").Value = 3517
$wsAll.Range("F19").Value = 233
